$d = $word.ActiveDocument

# The document contains two "<id>...</id>" blocks, each split across three
# runs: "<id>", the bare identifier text, and "</id>". Re-download bumped
# the identifiers from "p056r_a1"/"p056r_a2" to "p056r_1"/"p056r_2". Find the
# full "<id>...</id>" span (Word's Find happily matches across run
# boundaries) and overwrite its Range.Text in one shot — this collapses the
# three source runs into a single run that inherits the formatting of the
# first (and last) of them, i.e. the Courier-New/7f6000 "tag" styling, which
# matches the target markup.

$r1 = $d.Content.Duplicate
$found1 = $r1.Find.Execute("<id>p056r_a1</id>", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $r1.Text = "<id>p056r_1</id>"
}

$r2 = $d.Content.Duplicate
$found2 = $r2.Find.Execute("<id>p056r_a2</id>", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $r2.Text = "<id>p056r_2</id>"
}

Write-Output "id1 replaced: $found1"
Write-Output "id2 replaced: $found2"
